# The workbook stores test credentials on sheet "DataTest":
#   A1=Key        B1=username   C1=password
#   A2=UnitTest1  B2=student    C2=Password1234
# Update the stored password value (C2) to match the new test data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataTest")
$ws.Range("C2").Value = "Password123"
